# Add missing 2017 year data to the "CtIEPpUESoS" summary sheet.
#
# The sheet currently starts its year series at 2018 (column B). We insert
# a new column B for 2017, shifting the existing years (2018-2050) one
# column to the right (2019-2050 stay where a plain "insert column" would
# put them, ending at column AI instead of AH). The new 2017 column simply
# mirrors the (new) first data column C, matching how the other "roll
# forward" columns in each row reference the prior column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CtIEPpUESoS")

# Insert a new blank column before column B; this automatically shifts all
# existing formulas/values right by one column (B->C, C->D, ... AH->AI) and
# rewrites relative/absolute references (e.g. "=$B2" becomes "=$C2").
$ws.Columns("B").Insert()

# Header: 2017 for the newly inserted column.
$ws.Range("B1").Value = 2017

# Data rows: new column B simply equals the (new) first data column C, same
# pattern used elsewhere in these rows to carry a value forward across years.
$ws.Range("B2").Formula = "=C2"
$ws.Range("B3").Formula = "=C3"
$ws.Range("B4").Formula = "=C4"
$ws.Range("B5").Formula = "=C5"
$ws.Range("B6").Formula = "=C6"
$ws.Range("B7").Formula = "=C7"
$ws.Range("B8").Formula = "=C8"

# Match column A's width for the newly inserted column.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Make "CtIEPpUESoS" the active sheet/tab, with B2:B8 selected (mirrors the
# author reviewing the newly added 2017 figures).
$ws.Activate() | Out-Null
$ws.Range("B2:B8").Select() | Out-Null
